$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("distance from Dma50")

$values = @{
    2  = 10.2327
    3  = 7.4537
    4  = 6.0855
    5  = 5.2698
    6  = 5.1704
    7  = 4.9962
    8  = 4.5686
    9  = 4.4179
    10 = 3.9106
    11 = 3.5052
    12 = 3.415
    13 = 3.3895
    14 = 3.1123
    15 = 3.0775
    16 = 2.9969
    17 = 2.8412
    18 = 2.5675
    19 = 2.5356
    20 = 2.3676
    21 = 2.3194
    22 = 1.575
    24 = 1.3737
    25 = 1.19
    26 = 1.0891
    27 = 1.0727
    28 = 0.6366000000000001
    29 = 0.3541
    30 = -2.1236
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
